# Actualización automática 2025-08-14 15:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("C25").Value = 518.4
$wsGrupo.Range("L30").Value = 855.36
$wsGrupo.Range("M30").Value = 727.83
$wsGrupo.Range("D37").Value = 2747.52
$wsGrupo.Range("L57").Value = "4 de 55"
$wsGrupo.Range("M57").Value = "10 de 55"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F25").Value = 3049.78
$wsMensual.Range("F30").Value = 1583.19
$wsMensual.Range("F37").Value = 11710.93
$wsMensual.Range("F57").Value = 33979.48

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 - 240X120 PORCELANATO
$wsCumpl.Range("D2").Value = 518.4
$wsCumpl.Range("E2").Value = 9451.94304517915
$wsCumpl.Range("F2").Value = 0.05199419896095313

# Row 3 - 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 8363.51
$wsCumpl.Range("E3").Value = 19093.4976
$wsCumpl.Range("F3").Value = 0.3046038418257931

# Row 15 - PIEDRA SINTERIZADA
$wsCumpl.Range("D15").Value = 4452.57
$wsCumpl.Range("E15").Value = 9047.43
$wsCumpl.Range("F15").Value = 0.32982

# Row 16 - PORCELANATO
$wsCumpl.Range("D16").Value = 17699.02
$wsCumpl.Range("E16").Value = 38360.67999999999
$wsCumpl.Range("F16").Value = 0.3157173513236782

# Row 19 - TOTAL
$wsCumpl.Range("D19").Value = 33979.48
$wsCumpl.Range("E19").Value = 83460.21064517915
$wsCumpl.Range("F19").Value = 0.2893355714182038
